$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: APPROVAL_1 / APPROVAL_2 set to HOLD
$ws.Range("AI2").Value = "HOLD"
$ws.Range("AJ2").Value = "HOLD"

# ADJUSTMENT_AMOUNT (column AE) updates for pending rows
$ws.Range("AE3").Value = 150000
$ws.Range("AE4").Value = 132000
$ws.Range("AE5").Value = 14500
$ws.Range("AE6").Value = 1624.78
$ws.Range("AE7").Value = 1016.52
$ws.Range("AE8").Value = 13000
$ws.Range("AE9").Value = 79930
$ws.Range("AE10").Value = 560433
$ws.Range("AE13").Value = 80
$ws.Range("AE14").Value = 100
$ws.Range("AE15").Value = 250
$ws.Range("AE16").Value = 243
$ws.Range("AE17").Value = 23
$ws.Range("AE18").Value = 299
$ws.Range("AE22").Value = 1600
$ws.Range("AE24").Value = 724
$ws.Range("AE27").Value = 2500
$ws.Range("AE31").Value = 943
$ws.Range("AE32").Value = 1061
$ws.Range("AE35").Value = 1974
$ws.Range("AE36").Value = 2000
$ws.Range("AE37").Value = 850
$ws.Range("AE38").Value = 850
$ws.Range("AE39").Value = 30000
$ws.Range("AE40").Value = 14183
$ws.Range("AE41").Value = 14490
$ws.Range("AE42").Value = 15627.5
$ws.Range("AE43").Value = 15627.5
$ws.Range("AE44").Value = 13221
$ws.Range("AE45").Value = 16021
$ws.Range("AE46").Value = 18077.5
$ws.Range("AE47").Value = 11821
$ws.Range("AE48").Value = 13923
$ws.Range("AE49").Value = 15277
$ws.Range("AE50").Value = 13133
$ws.Range("AE51").Value = 1767.64
$ws.Range("AE52").Value = 13570
$ws.Range("AE53").Value = 16000
$ws.Range("AE54").Value = 342495
$ws.Range("AE55").Value = 5109
$ws.Range("AE56").Value = 396
$ws.Range("AE57").Value = 94
$ws.Range("AE58").Value = 1704
$ws.Range("AE59").Value = 94
$ws.Range("AE60").Value = 892
$ws.Range("AE61").Value = 94
$ws.Range("AE62").Value = 94
$ws.Range("AE63").Value = 94
$ws.Range("AE64").Value = 94
$ws.Range("AE65").Value = 1610
$ws.Range("AE66").Value = 94
$ws.Range("AE67").Value = 1502
$ws.Range("AE68").Value = 1274
$ws.Range("AE69").Value = 4236
$ws.Range("AE70").Value = 5516
$ws.Range("AE71").Value = 40
$ws.Range("AE72").Value = 23600
$ws.Range("AE73").Value = 1500
$ws.Range("AE74").Value = 1500
$ws.Range("AE76").Value = 2000

# Remove the last row (WGE 302 / Western Interior Designers, row 77)
$ws.Rows.Item(77).Delete()
